$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("area_pop_sum")

# Fix casing of the existing "Population" label -> "population"
$ws.Range("A3").Value = "population"

# Append a new "density" row
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = 3902.031942095055
